$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.724.96'
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").Value = '1.600.00'
$ws.Range("E3").Value = '  +0.12%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.06%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.513'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.06%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("E9").Value = '  +0.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.61'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.75%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0849'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.86%  '

$ws.Range("D12").Value = '1.823.63'
$ws.Range("E12").Value = '  +0.01%  '

$ws.Range("D13").Value = '1.601.43'
$ws.Range("E13").Value = '  +0.42%  '

$ws.Range("E14").Value = '  +0.80%  '

$ws.Range("E15").Value = '  +0.45%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.07'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.15%  '

$ws.Range("D17").Value = '0.0₃0741'
$ws.Range("E17").Value = '  -3.24%  '

$ws.Range("E18").Value = '  +0.06%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '208.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.17%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.84%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.66'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.38%  '

$ws.Range("E25").Value = '  -0.03%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.13%  '

$ws.Range("E27").Value = '  -0.70%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.34'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.17%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0508'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.15%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("E31").Value = '  +0.61%  '

$ws.Range("E32").Value = '  +0.64%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.27'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +18.62%  '

$ws.Range("D34").Value = '1.277.38'
$ws.Range("E34").Value = '  -0.57%  '

$ws.Range("E35").Value = '  +1.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.50'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.590'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.58%  '

$ws.Range("E38").Value = '  -1.24%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.824'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.12%  '

$ws.Range("E40").Value = '  +0.68%  '

$ws.Range("E41").Value = '  +0.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.778'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.83%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '62.70'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.68%  '

$ws.Range("D44").Value = '1.735.32'
$ws.Range("E44").Value = '  -0.04%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '90.48'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.69%  '

$ws.Range("E46").Value = '  +0.57%  '

$ws.Range("E47").Value = '  +1.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0512'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.95%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.17%  '

$ws.Range("E50").Value = '  +0.11%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.398'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.47%  '
